$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (2023-10-04 -> 2023-10-05, serial 45203 -> 45204) for every data row
# (rows 2 through 469) as part of an automatic daily update.
$firstRow = 2
$lastRow = 469

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
